# Apply updated dSF (column F) values for specific rows,
# as part of a "repull data, push all data, mean calculation" update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F7").Value  = -3
$ws.Range("F9").Value  = -1
$ws.Range("F12").Value = -4
$ws.Range("F15").Value = -6
$ws.Range("F16").Value = 2
$ws.Range("F17").Value = 1
$ws.Range("F18").Value = -5
